# Applies the commit "added season mapping O for one season":
#  1. Removes the empty "direct" (column E) placeholder cells from the
#     data rows of the existing six sheets (L, K1, K2, H1, H2, P).
#  2. Adds a brand-new worksheet named "O" at the end of the workbook
#     containing the consolidated "One Season" fare table (every
#     dest/booking_class/base_fare combination from the other six
#     sheets, re-tagged with season "O"), plus a remarks note.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: clear the stray empty "direct" cells (column E, rows 2..N)
# on each of the pre-existing sheets.
# ---------------------------------------------------------------------
$existingSheetNames = @("L", "K1", "K2", "H1", "H2", "P")
foreach ($name in $existingSheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $lastRow = $ws.Cells.Item(1, 1).End(-4121).Row
    if ($lastRow -ge 2) {
        $ws.Range("E2:E" + $lastRow).ClearContents()
    }
}

# ---------------------------------------------------------------------
# Step 2: add the new "O" worksheet after the last existing sheet ("P")
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "O"

# Header row
$headers = @("dest", "booking_class", "season", "base_fare", "direct", `
    "Unnamed: 6", "Unnamed: 7", "Unnamed: 8", "Unnamed: 9", "Unnamed: 10")
$colLetters = @("A","B","C","D","E","F","G","H","I","J")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Range($colLetters[$i] + "1")
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
}

# Data rows: one "O" season entry for every dest/booking_class/base_fare
# combination found across the other six sheets, ordered by base_fare.
$data = @(
    @("TPE","Y","O",2300),
    @("TPE","B","O",2400),
    @("TPE","M","O",2500),
    @("TPE","Q","O",2600),
    @("TPE","U","O",2700),
    @("TPE","H","O",2800),
    @("TPE","N","O",2900),
    @("TPE","A","O",3000),
    @("TPE","E","O",3100),
    @("SGN","Y","O",3200),
    @("SGN","B","O",3300),
    @("SGN","M","O",3400),
    @("SGN","Q","O",3500),
    @("SGN","U","O",3600),
    @("SGN","H","O",3700),
    @("SGN","N","O",3800),
    @("BKK","Y","O",3900),
    @("BKK","B","O",4000),
    @("BKK","M","O",4100),
    @("BKK","Q","O",4200),
    @("BKK","U","O",4300),
    @("BKK","H","O",4400),
    @("BKK","N","O",4500)
)

$rowIdx = 2
foreach ($rec in $data) {
    $ws.Range("A" + $rowIdx).Value = $rec[0]
    $ws.Range("B" + $rowIdx).Value = $rec[1]
    $ws.Range("C" + $rowIdx).Value = $rec[2]
    $ws.Range("D" + $rowIdx).Value = $rec[3]
    $rowIdx++
}

# Remarks note, attached to the first data row (column G)
$ws.Range("G2").Value = "Remarks:`n'O' for One Season"

Write-Host "Added sheet 'O' with" $data.Length "fare rows"
